$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (CellRef, NewValue) pairs. A $null value
# means the cell is cleared (matches the source diff removing that <c> node).
$updates = @{}
$updates["ALC"] = @(
    @("H17", 1561.875),
    @("J17", 1572.826),
    @("L17", 4718.478),
    @("N17", -5054.478),
    @("H40", 3873.5),
    @("J40", 4837.8),
    @("L40", 4837.8),
    @("N40", -5187.8),
    @("H42", 160.125),
    @("I42", 116.75),
    @("J42", 203.5),
    @("K42", 350.25),
    @("L42", 610.5),
    @("M42", -120.25),
    @("N42", -1070.5),
    @("H69", 7224.5),
    @("J69", 0),
    @("L69", 0),
    @("N69", $null),
    @("H72", 7224.5),
    @("J72", 0),
    @("L72", 0),
    @("N72", $null),
    @("H74", 6126),
    @("I74", 5377),
    @("K74", 5377),
    @("M74", -4441),
    @("H77", 6126),
    @("I77", 5377),
    @("K77", 26885),
    @("M77", -22205),
    @("H100", 89840.64999999999),
    @("I100", 200745.8),
    @("J100", 43630.168),
    @("K100", 200745.8),
    @("L100", 43630.168),
    @("M100", -200204.8),
    @("N100", -44712.168),
    @("H106", 6179429.5),
    @("I106", 6864941),
    @("K106", 6864941),
    @("M106", -6864310),
    @("H133", 94000),
    @("J133", 94000),
    @("L133", 94000),
    @("N133", -104120),
    @("H137", 8063.317),
    @("I137", 10377.552),
    @("J137", 2470.5833),
    @("K137", 31132.656),
    @("L137", 7411.749899999999),
    @("M137", -28582.656),
    @("N137", -12511.7499)
)

$updates["ARM"] = @(
    @("H44", 80000),
    @("J44", 80000),
    @("L44", 80000),
    @("N44", -80976),
    @("H132", 3346.45),
    @("I132", 1302.3077),
    @("K132", 3906.9231),
    @("M132", -1376.9231),
    @("H139", 202057.33),
    @("J139", 202057.33),
    @("L139", 202057.33),
    @("N139", -212337.33)
)

$updates["BSM"] = @(
    @("H20", 2780.625),
    @("I20", 1834.7858),
    @("J20", 4104.8),
    @("K20", 1834.7858),
    @("L20", 4104.8),
    @("M20", -1587.7858),
    @("N20", -4598.8),
    @("H35", 54999.6),
    @("J35", 54999.6),
    @("L35", 54999.6),
    @("N35", -55619.6),
    @("H105", 76964.336),
    @("I105", 101734.91),
    @("J105", 8845.25),
    @("K105", 101734.91),
    @("L105", 8845.25),
    @("M105", -99987.91),
    @("N105", -12339.25),
    @("H134", 6419.032),
    @("I134", 8115.706),
    @("K134", 24347.118),
    @("M134", -21812.118)
)

$updates["CRP"] = @(
    @("H31", 6958.6),
    @("I31", 7870.643),
    @("K31", 7870.643),
    @("M31", -7575.643),
    @("H34", 6958.6),
    @("I34", 7870.643),
    @("K34", 7870.643),
    @("M34", -7668.643),
    @("H58", 3398.3333),
    @("I58", 3818.0715),
    @("K58", 3818.0715),
    @("M58", -3615.0715),
    @("H105", 96239.41),
    @("I105", 124240.82),
    @("K105", 124240.82),
    @("M105", -122493.82),
    @("H132", 1383.037),
    @("I132", 1322.5),
    @("K132", 3967.5),
    @("M132", -1437.5),
    @("H136", 3398.3333),
    @("I136", 3818.0715),
    @("K136", 11454.2145),
    @("M136", -8904.2145)
)

$updates["CUL"] = @(
    @("H107", 1292.7894),
    @("J107", 1292.7894),
    @("L107", 3878.3682),
    @("N107", -7718.3682),
    @("H113", 800.55554),
    @("J113", 882.5),
    @("L113", 2647.5),
    @("N113", -6987.5),
    @("I114", 995.8),
    @("J114", 4245),
    @("K114", 2987.4),
    @("L114", 12735),
    @("M114", 266.6000000000004),
    @("N114", -19243),
    @("H115", 1239.4),
    @("J115", 0),
    @("L115", 0),
    @("N115", $null),
    @("H116", 4031.3333),
    @("I116", 2837.6),
    @("K116", 8512.799999999999),
    @("M116", -5070.799999999999)
)

$updates["GSM"] = @(
    @("H45", 29999.666),
    @("J45", 29999.666),
    @("L45", 29999.666),
    @("N45", -31117.666),
    @("H70", 14175.77),
    @("I70", 15262),
    @("J70", 13496.875),
    @("K70", 15262),
    @("L70", 13496.875),
    @("M70", -14992),
    @("N70", -14036.875),
    @("H73", 14175.77),
    @("I73", 15262),
    @("J73", 13496.875),
    @("K73", 15262),
    @("L73", 13496.875),
    @("M73", -14326),
    @("N73", -15368.875),
    @("H80", 18078.5),
    @("I80", 31799.2),
    @("K80", 31799.2),
    @("M80", -30801.2),
    @("H83", 18078.5),
    @("I83", 31799.2),
    @("K83", 158996),
    @("M83", -154004),
    @("H132", 5373.263),
    @("I132", 6020.1875),
    @("K132", 18060.5625),
    @("M132", -15530.5625)
)

$updates["LTW"] = @(
    @("H46", 2250.16),
    @("I46", 738.3333),
    @("K46", 738.3333),
    @("M46", -550.3333),
    @("H82", 3557.6),
    @("I82", 6220),
    @("J82", 1782.6666),
    @("K82", 6220),
    @("L82", 1782.6666),
    @("M82", -5859),
    @("N82", -2504.6666),
    @("H85", 3557.6),
    @("I85", 6220),
    @("J85", 1782.6666),
    @("K85", 6220),
    @("L85", 1782.6666),
    @("M85", -4972),
    @("N85", -4278.6666),
    @("H132", 483020.75),
    @("I132", 1065173.1),
    @("K132", 3195519.3),
    @("M132", -3192989.3),
    @("H136", 5098.485),
    @("I136", 3559.182),
    @("K136", 10677.546),
    @("M136", -8127.545999999998)
)

$updates["WVR"] = @(
    @("H81", 8517.786),
    @("I81", 11377.8),
    @("K81", 22755.6),
    @("M81", -21694.6),
    @("H82", 27000),
    @("J82", 27000),
    @("L82", 27000),
    @("N82", -27766),
    @("H84", 8517.786),
    @("I84", 11377.8),
    @("K84", 113778),
    @("M84", -108474),
    @("H85", 27000),
    @("J85", 27000),
    @("L85", 27000),
    @("N85", -29652)
)

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($pair in $updates[$sheetName]) {
        $cellRef = $pair[0]
        $newVal = $pair[1]
        if ($newVal -eq $null) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $newVal
        }
    }
}

Write-Output "Applied $($updates.Keys.Count) sheet updates"